# Update the "Forecast Comparison" sheet:
#  - Insert a new "Week_Start_Date" column after "Week" (new column B)
#  - Re-number the Week labels without the leading zero (W01 -> W1, etc.)
#  - Populate the new Week_Start_Date column with the week's start date (as text)
#  - Convert the is_holiday_week column to boolean-typed values

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new blank column at B; this shifts ASIN..is_holiday_week one column right
# (C..J) and also shifts the dimension/ranges automatically.
$ws.Columns("B:B").Insert()

# New header for the inserted column
$ws.Range("B1").Value = "Week_Start_Date"

# Week start dates (Sundays), one per data row (rows 2-17 => weeks 1-16)
$weekStartDates = @(
    "2025-01-05",
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20"
)

for ($i = 0; $i -lt $weekStartDates.Length; $i++) {
    $row = $i + 2
    $weekNum = $i + 1

    # Shorten the Week label (W01 -> W1, ... W16 stays W16)
    $ws.Range("A$row").Value = "W$weekNum"

    # Write the week start date as plain text (avoid Excel's automatic date
    # conversion by temporarily forcing a text number format, then reset the
    # cell style back to Normal so no lingering formatting remains).
    $dateCell = $ws.Range("B$row")
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $weekStartDates[$i]
    $dateCell.Style = "Normal"

    # is_holiday_week now lives in column J; store as a real boolean value
    $ws.Range("J$row").Value = $false
}

Write-Output "Updated Forecast Comparison sheet with Week_Start_Date column"
